$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New account-statement rows (worker order swapped: MARCO POLO HERRERA BERMEJO
# now listed first, then DANIEL CAMILO CASTRO REALES; periods listed in
# descending order 2304..2210; the 54000 reduced-value period is now the most
# recent period, 2304, for each worker).
$rows = @(
    @{ Row = 16; Doc = "1007048317"; Name = "MARCO POLO HERRERA BERMEJO"; Period = "2304"; Valor = 54000 },
    @{ Row = 17; Doc = "1007048317"; Name = "MARCO POLO HERRERA BERMEJO"; Period = "2303"; Valor = 60000 },
    @{ Row = 18; Doc = "1007048317"; Name = "MARCO POLO HERRERA BERMEJO"; Period = "2302"; Valor = 60000 },
    @{ Row = 19; Doc = "1007048317"; Name = "MARCO POLO HERRERA BERMEJO"; Period = "2301"; Valor = 60000 },
    @{ Row = 20; Doc = "1007048317"; Name = "MARCO POLO HERRERA BERMEJO"; Period = "2212"; Valor = 60000 },
    @{ Row = 21; Doc = "1007048317"; Name = "MARCO POLO HERRERA BERMEJO"; Period = "2211"; Valor = 60000 },
    @{ Row = 22; Doc = "1007048317"; Name = "MARCO POLO HERRERA BERMEJO"; Period = "2210"; Valor = 60000 },
    @{ Row = 23; Doc = "1143390341"; Name = "DANIEL CAMILO CASTRO REALES"; Period = "2304"; Valor = 54000 },
    @{ Row = 24; Doc = "1143390341"; Name = "DANIEL CAMILO CASTRO REALES"; Period = "2303"; Valor = 60000 },
    @{ Row = 25; Doc = "1143390341"; Name = "DANIEL CAMILO CASTRO REALES"; Period = "2302"; Valor = 60000 },
    @{ Row = 26; Doc = "1143390341"; Name = "DANIEL CAMILO CASTRO REALES"; Period = "2301"; Valor = 60000 },
    @{ Row = 27; Doc = "1143390341"; Name = "DANIEL CAMILO CASTRO REALES"; Period = "2212"; Valor = 60000 },
    @{ Row = 28; Doc = "1143390341"; Name = "DANIEL CAMILO CASTRO REALES"; Period = "2211"; Valor = 60000 },
    @{ Row = 29; Doc = "1143390341"; Name = "DANIEL CAMILO CASTRO REALES"; Period = "2210"; Valor = 60000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc      # column C: N° Doc Trabajador
    $ws.Cells.Item($r.Row, 4).Value = $r.Name     # column D: Nombre Trabajador
    $ws.Cells.Item($r.Row, 5).Value = $r.Period   # column E: Periodo Mora
    $ws.Cells.Item($r.Row, 6).Value = $r.Valor    # column F: Valor Mora
}
